$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = "66.778.85"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3: 'Ethereum'
$ws.Range("D3").Value = "3.187.49"
$ws.Range("E3").Value = "  -0.80%  "

# Row 4: 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5: 'BNB'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.08%  "

# Row 6: 'Solana'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.58%  "

# Row 7: 'USDC'
$ws.Range("E7").Value = "  +0.24%  "

# Row 8: 'LidoStakedEther'
$ws.Range("D8").Value = "3.166.13"
$ws.Range("E8").Value = "  -1.27%  "

# Row 9: 'XRP'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.91%  "

# Row 10: 'Dogecoin'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.94%  "

# Row 11: 'Toncoin'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.75%  "

# Row 12: 'Cardano'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.52%  "

# Row 13: 'ShibaInu'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.15%  "

# Row 14: 'Avalanche'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -8.11%  "

# Row 15: 'WrappedliquidstakedEther2.0'
$ws.Range("D15").Value = "3.707.77"
$ws.Range("E15").Value = "  -0.85%  "

# Row 16: 'WrappedBTC'
$ws.Range("D16").Value = "66.716.06"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17: 'WrappedEther'
$ws.Range("D17").Value = "3.200.76"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18: 'TRON'
$ws.Range("E18").Value = "  -1.56%  "

# Row 19: 'Polkadot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "

# Row 20: 'BitcoinCash'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "494.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.43%  "

# Row 21: 'Chainlink'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.98%  "

# Row 22: 'Polygon'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.04%  "

# Row 23: 'Uniswap'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.82%  "

# Row 24: 'Litecoin'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.79%  "

# Row 25: 'InternetComputer(DFINITY)'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.90%  "

# Row 26: 'Dai'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

# Row 27: 'PancakeSwap'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.76%  "

# Row 28: 'EthereumClassic'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "27.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.67%  "

# Row 29: 'ImmutableX'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.14%  "

# Row 30: 'RenderToken'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.38%  "

# Row 31: 'Mantle'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.54%  "

# Row 32: 'Stacks'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.57%  "

# Row 33: 'FirstDigitalUSD'
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.09%  "

# Row 34: 'Bittensor'
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "512.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.78%  "

# Row 35: 'OKB'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.41%  "

# Row 36: 'Filecoin'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.47%  "

# Row 37: 'NEARProtocol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.04%  "

# Row 38: 'VeChain'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0410"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.45%  "

# Row 39: 'Hedera'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0801"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.41%  "

# Row 40: 'Cosmos'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.46%  "

# Row 41: 'Kaspa'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.117"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.21%  "

# Row 42: 'Maker'
$ws.Range("D42").Value = "2.832.89"
$ws.Range("E42").Value = "  -1.11%  "

# Row 43: 'dogwifhat'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.79%  "

# Row 44: 'USDe'
$ws.Range("E44").Value = "  -0.07%  "

# Row 45: 'TheGraph'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.245"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.04%  "

# Row 46: 'Monero'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "120.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.64%  "

# Row 47: 'InjectiveProtocol'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.32%  "

# Row 48: 'Stellar'
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.108"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.13%  "

# Row 49: 'Fetch.AI'
$ws.Range("B49").Value = "Fetch.AI"
$ws.Range("C49").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.95%  "

# Row 50: 'PEPE'
$ws.Range("D50").Value = "0.0₃0514"
$ws.Range("E50").Value = "  -10.70%  "

# Row 51: 'Cronos'
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.136"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.57%  "
